$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feeds")
Write-Host $ws.Name
